$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = [double]45177

for ($r = 2; $r -le 295; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
